$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.386.63'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.86'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4670'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +1.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06549'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.04'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +8.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07929'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.29'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.874.52'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.152'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6746'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '282.14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.388.31'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.559'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.92%  '

$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("E20").Value = '  +1.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.123.21'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007293'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.26%  '

$ws.Range("E23").Value = '  +0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.206'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.295'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.09'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.937'
$ws.Range("D28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.355'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09687'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.434'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("E32").Value = '  -1.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.118'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04706'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.119'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7049'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.720'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01858'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.356'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -5.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.546'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.75'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("E42").Value = '  -0.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8473'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4189'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.58%  '

$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.72'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.220'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.248'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '930.47'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -6.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.13'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1132'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.03%  '
